$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has two "template"/placeholder rows at the top of the
# data block (row 2 = "Manual Posting" labels, row 3 = "(read only)" labels)
# before the real per-student rows start. Formatting the sheet for input
# means dropping those two placeholder rows so the real student data moves
# up and starts right under the header row.
$ws.Range("A2:A3").EntireRow.Delete()

# This file now represents the "02" section's course data export, so update
# the sheet/tab name to match.
$ws.Name = "COMP-101_FA24_02_course_data"
